$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: fill in the new week's coverage data (was just a bare index) ---
$ws.Range("A10").Value = 44081
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D10").Value = 162
$ws.Range("E10").Value = 78
$ws.Range("F10").Value = 1455
$ws.Range("G10").Value = 549
$ws.Range("H10").Value = 40
$ws.Range("I10").Value = 20
$ws.Range("J10").Value = 14
$ws.Range("K10").Value = 210
$ws.Range("L10").Value = 83

$ws.Range("N10").Formula = "=100*E10/D10"
$ws.Range("O10").Formula = "=100*G10/F10"
$ws.Range("P10").Formula = "=100*H10/D10"
$ws.Range("Q10").Formula = "=100*J10/I10"
$ws.Range("R10").Formula = "=100*L10/K10"

$ws.Range("N10:R10").NumberFormat = $ws.Range("N9:R9").NumberFormat

# --- Rows 11-41: the trailing placeholder index numbers (10..40) in column B
#     are no longer needed now that row 10 has real data, so clear them but
#     keep the cell (and its formatting) in place. ---
for ($row = 11; $row -le 41; $row++) {
    $ws.Cells.Item($row, 2).ClearContents()
}

# --- Update the active selection left on the sheet ---
$ws.Range("P23").Select()
